$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 96 (a new weekly price report), which
# pushes the existing rows 96..190 down to 97..191.
$ws.Rows("96").Insert()

# Populate the newly inserted row 96 with the new week's data. The
# "metadata" columns (A,B,C,E,F,G,H,I,J,K,L,Q,T) are constant across this
# product block, matching the rows around it.
$ws.Range("A96").Value = 5
$ws.Range("B96").Value = "Macroferia Regional de Talca"
$ws.Range("C96").Value = "Maule"
$ws.Range("D96").Value = 45167
$ws.Range("E96").Value = 7
$ws.Range("F96").Value = "Fruta"
$ws.Range("G96").Value = 100108
$ws.Range("H96").Value = "Tropicales y subtropicales"
$ws.Range("I96").Value = 100108002
$ws.Range("J96").Value = "Mango"
$ws.Range("K96").Value = "Sin especificar"
$ws.Range("L96").Value = "Primera"
$ws.Range("M96").Value = 248
$ws.Range("N96").Value = 9000
$ws.Range("O96").Value = 9000
$ws.Range("P96").Value = 9000
$ws.Range("Q96").Value = "$/bandeja 4 kilos"
$ws.Range("R96").Value = "Brasil"
$ws.Range("S96").Value = 2250
$ws.Range("T96").Value = 4
